$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated histogram bin counts (column B, rows 1-41) to match the recomputed data
$values = @(
    0,
    0,
    0,
    16,
    0,
    5,
    305,
    719,
    808,
    570,
    1080,
    817,
    680,
    275,
    353,
    696,
    1148,
    2042,
    2192,
    2637,
    3872,
    5790,
    6045,
    4763,
    3008,
    2435,
    2280,
    1910,
    1457,
    853,
    358,
    88,
    79,
    166,
    496,
    1677,
    3440,
    9064,
    18880,
    37808,
    287194
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $values[$i]
}

$excel.Calculate()
